# Fall 2022 Week 12 "day-after" inputs: fill in the M column (week of
# 11/15/2022) results for both tables, replacing the placeholder "A"
# (Available) values with the actual outcomes recorded the day after.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")

# --- Table 1 (rows 3-10) ---
$ws.Range("M3").Value = "L"
$ws.Range("M4").Value = "W"
$ws.Range("M5").Value = "NA"
$ws.Range("M6").Value = "L"
$ws.Range("M7").Value = "NA"
$ws.Range("M8").Value = "W"
$ws.Range("M9").Value = "DNP"
$ws.Range("M10").Value = "W"

# --- Table 2 (rows 15-22) ---
$ws.Range("M15").Value = "W"
$ws.Range("M16").Value = "W"
$ws.Range("M17").Value = "DNP"
$ws.Range("M18").Value = "NA"
$ws.Range("M19").Value = "DNP"
$ws.Range("M20").Value = "L"
$ws.Range("M21").Value = "W"
$ws.Range("M22").Value = "W"

# Match the saved view state captured after entering the data: the user
# scrolled down slightly and moved the frozen pane / selection over to
# column O, ending with V23 selected.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A4").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("O1").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("V23").Select()
